$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.625368237495422
$ws.Range("B1").Value = 1.914119839668274
$ws.Range("C1").Value = 2.016417980194092
$ws.Range("D1").Value = 2.36760425567627
$ws.Range("E1").Value = 3.216626405715942
